# Updated to fix panadas SettingWithCopy warning
#
# - Duplicate sheet "190824-00" to a new trailing sheet "190824-01"
#   (a new experiment run), widen it from 4 data columns (B:E) to
#   8 data columns (B:I), and clear out the per-seed columns that this
#   new run hasn't populated yet (only column B / seed S001 is filled
#   in), adjusting a handful of the remaining B-column config values.
# - Update the previously-active sheet's selection / tab-selected state
#   now that the new sheet is the active one.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("190824-00")

# Copy the source sheet to the very end of the workbook, then rename it.
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "190824-01"

# --- update the scalar config values that differ for this run ---
$ws3.Range("B3").Value = 640
$ws3.Range("B4").Value = 640
$ws3.Range("B12").Value = "640, 640"

# foodFile / mapFile are swapped relative to the source sheet
$ws3.Range("B8").Value = "foodMaps-04.png"
$ws3.Range("B9").Value = "foodMaps-00.png"

# --- row 2 ("completed"): blank out seeds S002:S004 down to 0, no style ---
$ws3.Range("C2:E2").Value = 0
$ws3.Range("C2:E2").Style = "Normal"

# --- rows 3:7 and 13 only ever had a single seed (column B); drop C:E ---
$ws3.Range("C3:E7").Clear()
$ws3.Range("C13:E13").Clear()

# --- rows 8:12 keep their formatting (style) in C:E, but are blanked ---
$ws3.Range("C8:E12").ClearContents()

# --- extend the run to 8 seed columns (F:I), matching C:E's blank style ---
$ws3.Range("C8:C12").Copy()
$ws3.Range("F8:I12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws3.Range("F1:I13").ColumnWidth = $ws3.Range("B1").ColumnWidth

# The previously-active sheet keeps a plain (non-tab-selected) view, but
# its remembered selection moves to E21.
$ws2 = $wb.Worksheets.Item("190824-00")
$ws2.Range("E21").Select() | Out-Null

# Finally, select the cell the new sheet should land on; this also marks
# the sheet as the active / tab-selected one and updates the workbook's
# activeTab.
$ws3.Range("D8").Select() | Out-Null
